# Auto-generated edit script: Add data for 2023-11-16
# Updates 2023 (column J) violent-crime-ytd figures across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 6766
$ws.Range('J3').Value = 7149
$ws.Range('J4').Value = 1557
$ws.Range('J5').Value = 560
$ws.Range('J6').Value = 9531
$ws.Range('J7').Value = 25563

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J4').Value = 115
$ws.Range('J7').Value = 741
$ws.Range('J8').Value = 1605
$ws.Range('J10').Value = 186
$ws.Range('J11').Value = 440
$ws.Range('J12').Value = 54
$ws.Range('J18').Value = 215
$ws.Range('J19').Value = 747
$ws.Range('J30').Value = 91
$ws.Range('J31').Value = 255
$ws.Range('J33').Value = 1154
$ws.Range('J36').Value = 348
$ws.Range('J37').Value = 791
$ws.Range('J39').Value = 17
$ws.Range('J41').Value = 178
$ws.Range('J42').Value = 1102
$ws.Range('J43').Value = 219
$ws.Range('J44').Value = 195
$ws.Range('J52').Value = 646
$ws.Range('J53').Value = 366
$ws.Range('J54').Value = 491
$ws.Range('J55').Value = 393
$ws.Range('J60').Value = 148
$ws.Range('J63').Value = 81
$ws.Range('J65').Value = 639
$ws.Range('J67').Value = 959
$ws.Range('J76').Value = 376
$ws.Range('J78').Value = 300
$ws.Range('J79').Value = 720
$ws.Range('J83').Value = 511
$ws.Range('J84').Value = 214
$ws.Range('J89').Value = 327
$ws.Range('J90').Value = 273
$ws.Range('J91').Value = 294
$ws.Range('J92').Value = 80
$ws.Range('J94').Value = 272
$ws.Range('J95').Value = 370
$ws.Range('J99').Value = 392
$ws.Range('J100').Value = 47
$ws.Range('J101').Value = 25563

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 226
$ws.Range('J7').Value = 741

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J6').Value = 201
$ws.Range('J7').Value = 440

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 99
$ws.Range('J6').Value = 96
$ws.Range('J7').Value = 327

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 153
$ws.Range('J6').Value = 275
$ws.Range('J7').Value = 646

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('J2').Value = 65
$ws.Range('J7').Value = 366

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 434
$ws.Range('J3').Value = 481
$ws.Range('J6').Value = 564
$ws.Range('J7').Value = 1605

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 151
$ws.Range('J7').Value = 511

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 268
$ws.Range('J3').Value = 384
$ws.Range('J6').Value = 401
$ws.Range('J7').Value = 1154

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J2').Value = 127
$ws.Range('J3').Value = 132
$ws.Range('J7').Value = 370

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 233
$ws.Range('J4').Value = 31
$ws.Range('J6').Value = 231
$ws.Range('J7').Value = 791

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 187
$ws.Range('J3').Value = 181
$ws.Range('J7').Value = 639

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J3').Value = 156
$ws.Range('J7').Value = 392

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('J6').Value = 23
$ws.Range('J7').Value = 91

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J2').Value = 89
$ws.Range('J7').Value = 255

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 355
$ws.Range('J6').Value = 268
$ws.Range('J7').Value = 959

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J2').Value = 65
$ws.Range('J7').Value = 214

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J3').Value = 102
$ws.Range('J4').Value = 36
$ws.Range('J7').Value = 491

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J2').Value = 181
$ws.Range('J6').Value = 289
$ws.Range('J7').Value = 747

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('J2').Value = 62
$ws.Range('J7').Value = 195

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J3').Value = 79
$ws.Range('J5').Value = 3
$ws.Range('J7').Value = 376

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J3').Value = 28
$ws.Range('J7').Value = 178

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J3').Value = 217
$ws.Range('J7').Value = 1102

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J6').Value = 105
$ws.Range('J7').Value = 186

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J3').Value = 94
$ws.Range('J7').Value = 300

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J2').Value = 76
$ws.Range('J6').Value = 219
$ws.Range('J7').Value = 393

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J6').Value = 74
$ws.Range('J7').Value = 294

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 200
$ws.Range('J6').Value = 215
$ws.Range('J7').Value = 720

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J3').Value = 43
$ws.Range('J7').Value = 215

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 112
$ws.Range('J7').Value = 348

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range('J3').Value = 12
$ws.Range('J7').Value = 47

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J6').Value = 145
$ws.Range('J7').Value = 272

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('J3').Value = 25
$ws.Range('J6').Value = 118

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range('J2').Value = 4
$ws.Range('J6').Value = 17

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('J6').Value = 23
$ws.Range('J7').Value = 80

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('J2').Value = 96
$ws.Range('J7').Value = 273

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('J2').Value = 52
$ws.Range('J7').Value = 148

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('J6').Value = 130
$ws.Range('J7').Value = 219

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('J2').Value = 36
$ws.Range('J6').Value = 50
$ws.Range('J7').Value = 115

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('J2').Value = 11
$ws.Range('J7').Value = 54
